$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated simulation metrics (Prey, Pred, Time(ms)) for generations 0-49 (rows 2-51)
$data = New-Object "object[,]" 50,3
$data[0,0] = 45
$data[0,1] = 8
$data[0,2] = 2854.290246963501
$data[1,0] = 49
$data[1,1] = 8
$data[1,2] = 3374.71079826355
$data[2,0] = 57
$data[2,1] = 9
$data[2,2] = 3491.153955459595
$data[3,0] = 70
$data[3,1] = 10
$data[3,2] = 3469.091176986694
$data[4,0] = 68
$data[4,1] = 11
$data[4,2] = 3477.364301681519
$data[5,0] = 73
$data[5,1] = 12
$data[5,2] = 3500.184535980225
$data[6,0] = 79
$data[6,1] = 8
$data[6,2] = 2777.769565582275
$data[7,0] = 78
$data[7,1] = 10
$data[7,2] = 2437.819242477417
$data[8,0] = 81
$data[8,1] = 16
$data[8,2] = 2425.56095123291
$data[9,0] = 79
$data[9,1] = 10
$data[9,2] = 2478.605508804321
$data[10,0] = 80
$data[10,1] = 13
$data[10,2] = 2467.833995819092
$data[11,0] = 82
$data[11,1] = 12
$data[11,2] = 2496.679067611694
$data[12,0] = 82
$data[12,1] = 16
$data[12,2] = 2461.416959762573
$data[13,0] = 84
$data[13,1] = 16
$data[13,2] = 2447.9820728302
$data[14,0] = 88
$data[14,1] = 14
$data[14,2] = 2436.282157897949
$data[15,0] = 88
$data[15,1] = 13
$data[15,2] = 2550.219774246216
$data[16,0] = 88
$data[16,1] = 12
$data[16,2] = 2463.279247283936
$data[17,0] = 85
$data[17,1] = 15
$data[17,2] = 2444.960117340088
$data[18,0] = 87
$data[18,1] = 15
$data[18,2] = 2455.043077468872
$data[19,0] = 88
$data[19,1] = 14
$data[19,2] = 2449.816703796387
$data[20,0] = 88
$data[20,1] = 15
$data[20,2] = 2461.83443069458
$data[21,0] = 89
$data[21,1] = 17
$data[21,2] = 2494.901180267334
$data[22,0] = 88
$data[22,1] = 17
$data[22,2] = 2500.065088272095
$data[23,0] = 88
$data[23,1] = 17
$data[23,2] = 2494.084358215332
$data[24,0] = 89
$data[24,1] = 17
$data[24,2] = 2460.748910903931
$data[25,0] = 88
$data[25,1] = 18
$data[25,2] = 2483.832836151123
$data[26,0] = 92
$data[26,1] = 17
$data[26,2] = 3773.667335510254
$data[27,0] = 92
$data[27,1] = 16
$data[27,2] = 2492.933034896851
$data[28,0] = 90
$data[28,1] = 18
$data[28,2] = 2464.917421340942
$data[29,0] = 93
$data[29,1] = 16
$data[29,2] = 2518.134117126465
$data[30,0] = 91
$data[30,1] = 18
$data[30,2] = 2447.538375854492
$data[31,0] = 91
$data[31,1] = 19
$data[31,2] = 2465.011119842529
$data[32,0] = 91
$data[32,1] = 18
$data[32,2] = 2468.804836273193
$data[33,0] = 91
$data[33,1] = 18
$data[33,2] = 2474.011898040771
$data[34,0] = 90
$data[34,1] = 18
$data[34,2] = 2438.58814239502
$data[35,0] = 92
$data[35,1] = 17
$data[35,2] = 2464.380264282227
$data[36,0] = 92
$data[36,1] = 16
$data[36,2] = 2481.318473815918
$data[37,0] = 91
$data[37,1] = 16
$data[37,2] = 2463.472127914429
$data[38,0] = 92
$data[38,1] = 16
$data[38,2] = 2474.099636077881
$data[39,0] = 92
$data[39,1] = 17
$data[39,2] = 2464.581727981567
$data[40,0] = 93
$data[40,1] = 17
$data[40,2] = 2509.955406188965
$data[41,0] = 93
$data[41,1] = 18
$data[41,2] = 2491.349458694458
$data[42,0] = 93
$data[42,1] = 17
$data[42,2] = 2484.604120254517
$data[43,0] = 93
$data[43,1] = 18
$data[43,2] = 2485.413074493408
$data[44,0] = 93
$data[44,1] = 19
$data[44,2] = 2466.079235076904
$data[45,0] = 93
$data[45,1] = 18
$data[45,2] = 2472.551584243774
$data[46,0] = 94
$data[46,1] = 18
$data[46,2] = 2490.83948135376
$data[47,0] = 94
$data[47,1] = 19
$data[47,2] = 2509.913206100464
$data[48,0] = 93
$data[48,1] = 19
$data[48,2] = 2462.328910827637
$data[49,0] = 94
$data[49,1] = 19
$data[49,2] = 2478.613615036011

$ws.Range("C2:E51").Value = $data

# Remove trailing generations 50-56 (rows 52-58), which no longer exist in the refactored simulation output
$ws.Range("A52:G58").Delete()
